$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data values (price & 1h volume change)
# Cells whose new value is purely numeric-looking need to be protected with a
# temporary Text number format so Excel stores them as exact text (matching the
# inlineStr cells in the source file) instead of converting them to floating point
# numbers (which would introduce binary rounding noise like "557.74000000000001").

$ws.Range('D2').Value = '68.391.74'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = '2.453.94'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.53%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.509'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.51%  '
$ws.Range('D9').Value = '2.452.86'
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('D10').Value = '0.157'
$ws.Range('E10').Value = '  +6.90%  '
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.81'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.15%  '
$ws.Range('B13').Value = 'Cardano'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.327'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.12%  '
$ws.Range('D14').Value = '68.352.31'
$ws.Range('E14').Value = '  +0.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000169'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.33'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '10.47'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '336.68'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.45%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('B22').Value = 'SuiNetwork'
$ws.Range('C22').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.87'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.96%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.59'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.65'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.10'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.88%  '
$ws.Range('D26').Value = '0.0₃0815'
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.48%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '422.27'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.11%  '
$ws.Range('E30').Value = '  +2.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.61'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('D32').Value = '161.03'
$ws.Range('E32').Value = '  +2.40%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = '17.75'
$ws.Range('E35').Value = '  +0.72%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '0.105'
$ws.Range('E36').Value = '  -1.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.295'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.85%  '
$ws.Range('E39').Value = '  +1.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.06'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.36'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '129.49'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0720'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.480'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.561'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0918'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.07%  '
$ws.Range('E48').Value = '  +1.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.36'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.67'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.72%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.87'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.79%  '
